# Automatic update of files.
# Bump the "Förändrad" (Changed) date column (C) by one day for all data
# rows (2 through 28): 45553 -> 45554 (serial date values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = 45554
}
